$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.346.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.368.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +6.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.31'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.49%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +26.33%  '
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.02'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +18.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +17.97%  '
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.720.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '17.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.72%  '
$ws.Range("E17").Value = '  +7.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.371.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '44.346.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("E20").Value = '  +3.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '78.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.59%  '
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.85%  '
$ws.Range("E28").Value = '  +3.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.58%  '
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("E33").Value = '  +6.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0757'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.45%  '
$ws.Range("E36").Value = '  +5.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.22%  '
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0278'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '19.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.34%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.197'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +17.38%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.101'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.69%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("E47").Value = '  +3.34%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.472.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
